# "actual final session logs"
#
# A new mVR session record was logged for animal 146 (session 135,
# "switchL20") right after its existing block of rows. In the saved
# worksheet this shows up as two new rows being inserted just above the
# animal-152 block (row 98 gets the new data, row 99 stays blank to
# preserve the usual one-blank-row separator between animal blocks), and
# every row from the old row 99 onward shifting down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows right before the old row 99 (push everything else
# down by two rows).
$ws.Rows("98:99").Insert()

# Row 98 carries the new record; row 99 is left blank, matching the
# existing blank-row separator convention between animal blocks.
$ws.Range("B98").Value = 146
$ws.Range("C98").Value = 135
$ws.Range("D98").Value = "switchL20"

# Restore the on-screen scroll/selection state as it was when the sheet
# was last saved.
$ws.Range("G90").Select()
$excel.ActiveWindow.ScrollRow = 74
$excel.ActiveWindow.Height = 14240
